$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J, matching style of existing headers (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) per row 2..18
$data = @(
    @(2, 2, 6),
    @(3, 9, 9),
    @(4, 3, 6),
    @(5, 3, 6),
    @(6, 9, 9),
    @(7, 3, 7),
    @(8, 7, 9),
    @(9, 3, 5),
    @(10, 1, 4),
    @(11, 1, 5),
    @(12, 1, 5),
    @(13, 1, 4),
    @(14, 1, 4),
    @(15, 3, 5),
    @(16, 3, 5),
    @(17, 6, 7),
    @(18, 8, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
